# Populate has_contrat with contrat <get all contrats>
#
# Moves the "real madrid fc" record up into row 2 (replacing "Ahmed Test"),
# updates its contract number and amounts, turns the old row 3 into the
# trailing blank/total row, and removes the old row 4 entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: becomes the "real madrid fc" record with an updated contract
# number (no "/AV1" suffix) and updated amounts.
$ws.Range("A2").Value = "real madrid fc"
$ws.Range("B2").Value = "110384"
$ws.Range("C2").Value = "114298485748578394873948"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "bmce"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "108/ANSYSFYSN01"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 40000
$ws.Range("J2").Value = 1800
$ws.Range("K2").Value = 38200

# Row 3: becomes the blank/total row with the same totals as row 2.
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("I3").Value = 40000
$ws.Range("J3").Value = 1800
$ws.Range("K3").Value = 38200

# Row 4: delete entirely, shifting dimension/used-range back down to K3.
$ws.Rows("4").Delete()
